$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7..129 down to 8..130
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the new record's data.
# Columns A,B,C,E,F,G,H,I,Q,R are identical to the rest of the dataset / unchanged,
# so copy them from the row directly below (old row 7, now row 8) to stay consistent,
# then set the changed fields (D, J, K, L, M, N, O, P).
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(7, $col).Value2 = $ws.Cells.Item(8, $col).Value2
}

# Fecha (D7) = 2021-10-27 -> Excel serial date 44496
$ws.Cells.Item(7, 4).Value2 = 44496
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat

$ws.Cells.Item(7, 10).Value2 = 150
$ws.Cells.Item(7, 11).Value2 = 900
$ws.Cells.Item(7, 12).Value2 = 900
$ws.Cells.Item(7, 13).Value2 = 900
$ws.Cells.Item(7, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(7, 15).Value2 = "Perú"
$ws.Cells.Item(7, 16).Value2 = 900
